$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5230.8096
$ws.Range("J17").Value = 5486.85
$ws.Range("L17").Value = 16460.55
$ws.Range("N17").Value = -16796.55
$ws.Range("H40").Value = 1770.3889
$ws.Range("J40").Value = 1294.9048
$ws.Range("L40").Value = 1294.9048
$ws.Range("N40").Value = -1644.9048
$ws.Range("H55").Value = 1166.6666
$ws.Range("J55").Value = 1166.6666
$ws.Range("L55").Value = 1166.6666
$ws.Range("N55").Value = -1594.6666
$ws.Range("H116").Value = 2023.9
$ws.Range("I116").Value = 1817.5
$ws.Range("J116").Value = 2333.5
$ws.Range("K116").Value = 1817.5
$ws.Range("L116").Value = 2333.5
$ws.Range("M116").Value = 1624.5
$ws.Range("N116").Value = -9217.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 66498.75
$ws.Range("J55").Value = 66498.75
$ws.Range("L55").Value = 66498.75
$ws.Range("N55").Value = -67128.75
$ws.Range("H102").Value = 4118615
$ws.Range("I102").Value = 4632892
$ws.Range("J102").Value = 4400
$ws.Range("K102").Value = 4632892
$ws.Range("L102").Value = 4400
$ws.Range("M102").Value = -4631270
$ws.Range("N102").Value = -7644
$ws.Range("H122").Value = 7354792.5
$ws.Range("I122").Value = 1797
$ws.Range("J122").Value = 31252028
$ws.Range("K122").Value = 5391
$ws.Range("L122").Value = 93756084
$ws.Range("M122").Value = -2941
$ws.Range("N122").Value = -93760984

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 246363.64
$ws.Range("I4").Value = 5000000
$ws.Range("J4").Value = 20000
$ws.Range("K4").Value = 5000000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = -4999888
$ws.Range("N4").Value = -20224
$ws.Range("H7").Value = 142.26666
$ws.Range("J7").Value = 224
$ws.Range("L7").Value = 224
$ws.Range("N7").Value = -450
$ws.Range("H31").Value = 5639.029
$ws.Range("I31").Value = 4881.8276
$ws.Range("J31").Value = 9298.833000000001
$ws.Range("K31").Value = 4881.8276
$ws.Range("L31").Value = 9298.833000000001
$ws.Range("M31").Value = -4586.8276
$ws.Range("N31").Value = -9888.833000000001
$ws.Range("H34").Value = 5639.029
$ws.Range("I34").Value = 4881.8276
$ws.Range("J34").Value = 9298.833000000001
$ws.Range("K34").Value = 4881.8276
$ws.Range("L34").Value = 9298.833000000001
$ws.Range("M34").Value = -4679.8276
$ws.Range("N34").Value = -9702.833000000001
$ws.Range("H107").Value = 1368.6666
$ws.Range("I107").Value = 2154
$ws.Range("J107").Value = 845.1111
$ws.Range("K107").Value = 2154
$ws.Range("L107").Value = 845.1111
$ws.Range("M107").Value = -234
$ws.Range("N107").Value = -4685.1111
$ws.Range("H122").Value = 6792.7
$ws.Range("I122").Value = 3317.423
$ws.Range("J122").Value = 29382
$ws.Range("K122").Value = 9952.269
$ws.Range("L122").Value = 88146
$ws.Range("M122").Value = -7502.269
$ws.Range("N122").Value = -93046
$ws.Range("H132").Value = 10706.267
$ws.Range("I132").Value = 14058.2
$ws.Range("J132").Value = 4002.4
$ws.Range("K132").Value = 42174.60000000001
$ws.Range("L132").Value = 12007.2
$ws.Range("M132").Value = -39644.60000000001
$ws.Range("N132").Value = -17067.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4303.0713
$ws.Range("J81").Value = 4395.25
$ws.Range("L81").Value = 13185.75
$ws.Range("N81").Value = -15431.75
$ws.Range("H84").Value = 4303.0713
$ws.Range("J84").Value = 4395.25
$ws.Range("L84").Value = 39557.25
$ws.Range("N84").Value = -50789.25
$ws.Range("H131").Value = 36190.332
$ws.Range("I131").Value = 1230.0834
$ws.Range("J131").Value = 64158.535
$ws.Range("K131").Value = 3690.2502
$ws.Range("L131").Value = 192475.605
$ws.Range("M131").Value = 1349.7498
$ws.Range("N131").Value = -202555.605

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9681.25
$ws.Range("I5").Value = 4990
$ws.Range("J5").Value = 17500
$ws.Range("K5").Value = 4990
$ws.Range("L5").Value = 17500
$ws.Range("M5").Value = -4878
$ws.Range("N5").Value = -17724
$ws.Range("H113").Value = 1786.16
$ws.Range("I113").Value = 1740.1904
$ws.Range("J113").Value = 2027.5
$ws.Range("K113").Value = 1740.1904
$ws.Range("L113").Value = 2027.5
$ws.Range("M113").Value = 429.8096
$ws.Range("N113").Value = -6367.5
$ws.Range("H122").Value = 6240.1763
$ws.Range("I122").Value = 8797.875
$ws.Range("K122").Value = 26393.625
$ws.Range("M122").Value = -23943.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 507
$ws.Range("I22").Value = 524.75
$ws.Range("J22").Value = 489.25
$ws.Range("K22").Value = 524.75
$ws.Range("L22").Value = 489.25
$ws.Range("M22").Value = -229.75
$ws.Range("N22").Value = -1079.25
$ws.Range("H27").Value = 507
$ws.Range("I27").Value = 524.75
$ws.Range("J27").Value = 489.25
$ws.Range("K27").Value = 524.75
$ws.Range("L27").Value = 489.25
$ws.Range("M27").Value = -417.75
$ws.Range("N27").Value = -703.25
$ws.Range("H46").Value = 1675
$ws.Range("I46").Value = 1675
$ws.Range("K46").Value = 1675
$ws.Range("M46").Value = -1487
$ws.Range("H55").Value = 315.33334
$ws.Range("I55").Value = 258.33334
$ws.Range("J55").Value = 372.33334
$ws.Range("K55").Value = 258.33334
$ws.Range("L55").Value = 372.33334
$ws.Range("M55").Value = -85.33334000000002
$ws.Range("N55").Value = -718.33334
$ws.Range("H93").Value = 1244.7778
$ws.Range("I93").Value = 967.1667
$ws.Range("J93").Value = 1800
$ws.Range("K93").Value = 967.1667
$ws.Range("L93").Value = 1800
$ws.Range("M93").Value = 280.8333
$ws.Range("N93").Value = -4296
$ws.Range("H122").Value = 6961.2
$ws.Range("I122").Value = 6345.8335
$ws.Range("J122").Value = 7884.25
$ws.Range("K122").Value = 19037.5005
$ws.Range("L122").Value = 23652.75
$ws.Range("M122").Value = -16587.5005
$ws.Range("N122").Value = -28552.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2490.6
$ws.Range("J107").Value = 3731.6667
$ws.Range("L107").Value = 11195.0001
$ws.Range("N107").Value = -15035.0001
